$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.043318213101276
$ws.Range("D2").Value = 1.045962672907326
$ws.Range("E2").Value = 1.041246706124911
$ws.Range("F2").Value = 1.049395666267763
$ws.Range("I2").Value = 1.036788601033545
$ws.Range("J2").Value = 1.048389304924778
$ws.Range("K2").Value = 1.048729441353454
$ws.Range("L2").Value = 1.044026763069418
$ws.Range("M2").Value = 1.05215284279171
$ws.Range("N2").Value = 1.049878137989705
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.045025910918264
$ws.Range("D3").Value = 1.047585557079674
$ws.Range("E3").Value = 1.042725417895944
$ws.Range("F3").Value = 1.051158425463515
$ws.Range("I3").Value = 1.037179307024748
$ws.Range("J3").Value = 1.049740042365018
$ws.Range("K3").Value = 1.05016176319882
$ws.Range("L3").Value = 1.045314309572809
$ws.Range("M3").Value = 1.053725386692558
$ws.Range("N3").Value = 1.051230793632042
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.046127905750465
$ws.Range("D4").Value = 1.048633018192531
$ws.Range("E4").Value = 1.043679734442459
$ws.Range("F4").Value = 1.052296303322753
$ws.Range("I4").Value = 1.037429441038133
$ws.Range("J4").Value = 1.050610776970278
$ws.Range("K4").Value = 1.051085449548366
$ws.Range("L4").Value = 1.046144431768193
$ws.Range("M4").Value = 1.054739743071068
$ws.Range("N4").Value = 1.052102764780316
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.046590481010069
$ws.Range("D5").Value = 1.049072748844756
$ws.Range("E5").Value = 1.044080340274463
$ws.Range("F5").Value = 1.052774024865017
$ws.Range("I5").Value = 1.037533960122873
$ws.Range("J5").Value = 1.050976059627987
$ws.Range("K5").Value = 1.051473032176704
$ws.Range("L5").Value = 1.046492706318832
$ws.Range("M5").Value = 1.055165429765128
$ws.Range("N5").Value = 1.052468566181298
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.046668108613301
$ws.Range("D6").Value = 1.049146545464382
$ws.Range("E6").Value = 1.044147569566487
$ws.Range("F6").Value = 1.052854199179081
$ws.Range("I6").Value = 1.037551472085828
$ws.Range("J6").Value = 1.051037347135061
$ws.Range("K6").Value = 1.051538066250931
$ws.Range("L6").Value = 1.046551141898685
$ws.Range("M6").Value = 1.055236860923855
$ws.Range("N6").Value = 1.052529940723662
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.046134089447216
$ws.Range("D7").Value = 1.048638896320891
$ws.Range("E7").Value = 1.043685089654768
$ws.Range("F7").Value = 1.052302689164817
$ws.Range("I7").Value = 1.0374308401256
$ws.Range("J7").Value = 1.050615660920867
$ws.Range("K7").Value = 1.051090631318141
$ws.Range("L7").Value = 1.046149088199942
$ws.Range("M7").Value = 1.054745434040737
$ws.Range("N7").Value = 1.052107655666676
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.043895967659783
$ws.Range("D8").Value = 1.04651169163739
$ws.Range("E8").Value = 1.041746970286806
$ws.Range("F8").Value = 1.049991976413026
$ws.Range("I8").Value = 1.036921199100896
$ws.Range("J8").Value = 1.048846479804285
$ws.Range("K8").Value = 1.049214154975895
$ws.Range("L8").Value = 1.044462524566288
$ws.Range("M8").Value = 1.052684959051752
$ws.Range("N8").Value = 1.050335962109979
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.039928468516762
$ws.Range("D9").Value = 1.042742365438468
$ws.Range("E9").Value = 1.038312010192486
$ws.Range("F9").Value = 1.045898508462701
$ws.Range("I9").Value = 1.036002445126845
$ws.Range("J9").Value = 1.045703283961065
$ws.Range("K9").Value = 1.045883104492734
$ws.Range("L9").Value = 1.041467064295643
$ws.Range("M9").Value = 1.04902913057492
$ws.Range("N9").Value = 1.047188302568133
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.037266635363784
$ws.Range("D10").Value = 1.040214576133379
$ws.Range("E10").Value = 1.036008013235276
$ws.Range("F10").Value = 1.043153998043618
$ws.Range("I10").Value = 1.035375764061804
$ws.Range("J10").Value = 1.043589813655547
$ws.Range("K10").Value = 1.043645183799177
$ws.Range("L10").Value = 1.039453583179566
$ws.Range("M10").Value = 1.046574229453325
$ws.Range("N10").Value = 1.045071830892363
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.036109829505709
$ws.Range("D11").Value = 1.039116291943823
$ws.Range("E11").Value = 1.035006863357959
$ws.Range("F11").Value = 1.041961699405263
$ws.Range("I11").Value = 1.035100982775788
$ws.Range("J11").Value = 1.042670219635152
$ws.Range("K11").Value = 1.042671881201568
$ws.Range("L11").Value = 1.038577652545029
$ws.Range("M11").Value = 1.045506841745804
$ws.Range("N11").Value = 1.04415093094298
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.035679488915992
$ws.Range("D12").Value = 1.038707762866084
$ws.Range("E12").Value = 1.034634450737414
$ws.Range("F12").Value = 1.041518221222452
$ws.Range("I12").Value = 1.034998397194114
$ws.Range("J12").Value = 1.042327958775622
$ws.Range("K12").Value = 1.042309696947719
$ws.Range("L12").Value = 1.038251666403589
$ws.Range("M12").Value = 1.045109688609546
$ws.Range("N12").Value = 1.043808184033772
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.035771828204438
$ws.Range("D13").Value = 1.038795420134084
$ws.Range("E13").Value = 1.034714359224552
$ws.Range("F13").Value = 1.041613376489314
$ws.Range("I13").Value = 1.035020425744527
$ws.Range("J13").Value = 1.042401406020184
$ws.Range("K13").Value = 1.042387416644239
$ws.Range("L13").Value = 1.038321620106967
$ws.Range("M13").Value = 1.045194910250992
$ws.Range("N13").Value = 1.043881735581845
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.036074270812478
$ws.Range("D14").Value = 1.039082534677328
$ws.Range("E14").Value = 1.034976090735163
$ws.Range("F14").Value = 1.041925053816688
$ws.Range("I14").Value = 1.035092513648697
$ws.Range("J14").Value = 1.042641942257135
$ws.Range("K14").Value = 1.042641956415089
$ws.Range("L14").Value = 1.038550719288185
$ws.Range("M14").Value = 1.045474026869214
$ws.Range("N14").Value = 1.044122613407843
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.036260528821849
$ws.Range("D15").Value = 1.039259358416743
$ws.Range("E15").Value = 1.035137280033871
$ws.Range("F15").Value = 1.042117007734
$ws.Range("I15").Value = 1.03513686038223
$ws.Range("J15").Value = 1.042790053670284
$ws.Range("K15").Value = 1.04279869934381
$ws.Range("L15").Value = 1.038691791429809
$ws.Range("M15").Value = 1.045645909504411
$ws.Range("N15").Value = 1.044270935156187
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.037343318414627
$ws.Range("D16").Value = 1.040287385562145
$ws.Range("E16").Value = 1.036074381192845
$ws.Range("F16").Value = 1.04323304309173
$ws.Range("I16").Value = 1.035393927810033
$ws.Range("J16").Value = 1.043650749144492
$ws.Range("K16").Value = 1.043709687477887
$ws.Range("L16").Value = 1.039511628698785
$ws.Range("M16").Value = 1.046644974336259
$ws.Range("N16").Value = 1.045132852916692
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.038021382585736
$ws.Range("D17").Value = 1.040931228469308
$ws.Range("E17").Value = 1.036661251660239
$ws.Range("F17").Value = 1.043932043566626
$ws.Range("I17").Value = 1.035554259051076
$ws.Range("J17").Value = 1.04418943977703
$ws.Range("K17").Value = 1.044279972960689
$ws.Range("L17").Value = 1.040024789121714
$ws.Range("M17").Value = 1.04727047164136
$ws.Range("N17").Value = 1.045672308551716
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.038416480543237
$ws.Range("D18").Value = 1.04130641190859
$ws.Range("E18").Value = 1.03700322618195
$ws.Range("F18").Value = 1.044339382947308
$ws.Range("I18").Value = 1.03564744751051
$ws.Range("J18").Value = 1.044503220786124
$ws.Range("K18").Value = 1.044612200046359
$ws.Range("L18").Value = 1.040323714173226
$ws.Range("M18").Value = 1.047634890315605
$ws.Range("N18").Value = 1.045986535165835
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.038551130520803
$ws.Range("D19").Value = 1.041434279324186
$ws.Range("E19").Value = 1.037119773917011
$ws.Range("F19").Value = 1.044478211963448
$ws.Range("I19").Value = 1.035679166540242
$ws.Range("J19").Value = 1.044610139820913
$ws.Range("K19").Value = 1.044725411674019
$ws.Range("L19").Value = 1.040425573681856
$ws.Range("M19").Value = 1.047759076370774
$ws.Range("N19").Value = 1.046093606037917
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.037948674781229
$ws.Range("D20").Value = 1.040862187493874
$ws.Range("E20").Value = 1.036598320985246
$ws.Range("F20").Value = 1.043857086430821
$ws.Range("I20").Value = 1.035537091200918
$ws.Range("J20").Value = 1.044131687757814
$ws.Range("K20").Value = 1.044218829324127
$ws.Range("L20").Value = 1.039969772585509
$ws.Range("M20").Value = 1.047203405607988
$ws.Range("N20").Value = 1.045614474518009
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.035985227131301
$ws.Range("D21").Value = 1.03899800267242
$ws.Range("E21").Value = 1.034899032380141
$ws.Range("F21").Value = 1.04183328947183
$ws.Range("I21").Value = 1.035071299942006
$ws.Range("J21").Value = 1.042571129293765
$ws.Range("K21").Value = 1.042567019036881
$ws.Range("L21").Value = 1.03848327269344
$ws.Range("M21").Value = 1.045391852829451
$ws.Range("N21").Value = 1.044051699881942
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.034746951767993
$ws.Range("D22").Value = 1.037822566633807
$ws.Range("E22").Value = 1.033827484915468
$ws.Range("F22").Value = 1.040557336216914
$ws.Range("I22").Value = 1.034775429790635
$ws.Range("J22").Value = 1.041585987064986
$ws.Range("K22").Value = 1.041524655404572
$ws.Range("L22").Value = 1.037545019641464
$ws.Range("M22").Value = 1.044248927560561
$ws.Range("N22").Value = 1.043065158638202
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.03540374919986
$ws.Range("D23").Value = 1.03844601047015
$ws.Range("E23").Value = 1.034395834779091
$ws.Range("F23").Value = 1.041234082228255
$ws.Range("I23").Value = 1.03493256307302
$ws.Range("J23").Value = 1.042108609587797
$ws.Range("K23").Value = 1.042077597837236
$ws.Range("L23").Value = 1.038042754125237
$ws.Range("M23").Value = 1.044855192085349
$ws.Range("N23").Value = 1.04358852334494
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.037981529560808
$ws.Range("D24").Value = 1.040893385248661
$ws.Range("E24").Value = 1.036626757689905
$ws.Range("F24").Value = 1.043890957494845
$ws.Range("I24").Value = 1.035544849633181
$ws.Range("J24").Value = 1.04415778473208
$ws.Range("K24").Value = 1.044246458770714
$ws.Range("L24").Value = 1.03999463340345
$ws.Range("M24").Value = 1.047233711169441
$ws.Range("N24").Value = 1.045640608552971
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.040957059822397
$ws.Range("D25").Value = 1.04371939159953
$ws.Range("E25").Value = 1.039202446152456
$ws.Range("F25").Value = 1.046959435784122
$ws.Range("I25").Value = 1.036242444834491
$ws.Range("J25").Value = 1.046518994002345
$ws.Range("K25").Value = 1.046747238722924
$ws.Range("L25").Value = 1.042244320767399
$ws.Range("M25").Value = 1.049977302236462
$ws.Range("N25").Value = 1.048005171011235
